$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 37.11928833333334
$ws.Range("H2").Value = 111.357865
$ws.Range("I2").Value = 0.08973251933053689
$ws.Range("J2").Value = 0.08973251933053687
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 62.351132
$ws.Range("N2").Value = 187.053396
$ws.Range("O2").Value = 0.9620005273240912
$ws.Range("P2").Value = 0.9620005273240912
$ws.Range("Q2").Value = 2314.429646617727
$ws.Range("R2").Value = 20829.86681955954
$ws.Range("S2").Value = 0.0863227309140957
$ws.Range("T2").Value = 0.08632273091409569

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.11928833333334
$ws.Range("H3").Value = 111.357865
$ws.Range("I3").Value = 0.08973251933053689
$ws.Range("J3").Value = 0.08973251933053687
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8562703333333334
$ws.Range("N3").Value = 2.568811
$ws.Range("O3").Value = 0.01321118776478095
$ws.Range("P3").Value = 0.01321118776478095
$ws.Range("Q3").Value = 31.78414539427945
$ws.Range("R3").Value = 286.057308548515
$ws.Range("S3").Value = 0.001185473161482559
$ws.Range("T3").Value = 0.001185473161482559

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 37.11928833333334
$ws.Range("H4").Value = 111.357865
$ws.Range("I4").Value = 0.08973251933053689
$ws.Range("J4").Value = 0.08973251933053687
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.606628666666667
$ws.Range("N4").Value = 4.819886
$ws.Range("O4").Value = 0.02478828491112776
$ws.Range("P4").Value = 0.02478828491112776
$ws.Range("Q4").Value = 59.63691272259889
$ws.Range("R4").Value = 536.7322145033901
$ws.Range("S4").Value = 0.002224315254958627
$ws.Range("T4").Value = 0.002224315254958627

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 168.3278913333333
$ws.Range("H5").Value = 504.983674
$ws.Range("I5").Value = 0.406917439453518
$ws.Range("J5").Value = 0.4069174394535179
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 62.351132
$ws.Range("N5").Value = 187.053396
$ws.Range("O5").Value = 0.9620005273240912
$ws.Range("P5").Value = 0.9620005273240912
$ws.Range("Q5").Value = 10495.43457180632
$ws.Range("R5").Value = 94458.9111462569
$ws.Range("S5").Value = 0.3914547913316533
$ws.Range("T5").Value = 0.3914547913316532

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 168.3278913333333
$ws.Range("H6").Value = 504.983674
$ws.Range("I6").Value = 0.406917439453518
$ws.Range("J6").Value = 0.4069174394535179
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8562703333333334
$ws.Range("N6").Value = 2.568811
$ws.Range("O6").Value = 0.01321118776478095
$ws.Range("P6").Value = 0.01321118776478095
$ws.Range("Q6").Value = 144.1341796212904
$ws.Range("R6").Value = 1297.207616591614
$ws.Range("S6").Value = 0.005375862697384311
$ws.Range("T6").Value = 0.00537586269738431

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 168.3278913333333
$ws.Range("H7").Value = 504.983674
$ws.Range("I7").Value = 0.406917439453518
$ws.Range("J7").Value = 0.4069174394535179
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.606628666666667
$ws.Range("N7").Value = 4.819886
$ws.Range("O7").Value = 0.02478828491112776
$ws.Range("P7").Value = 0.02478828491112776
$ws.Range("Q7").Value = 270.4404156156849
$ws.Range("R7").Value = 2433.963740541164
$ws.Range("S7").Value = 0.01008678542448038
$ws.Range("T7").Value = 0.01008678542448038

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 208.2187756666666
$ws.Range("H8").Value = 624.6563269999999
$ws.Range("I8").Value = 0.5033500412159452
$ws.Range("J8").Value = 0.5033500412159452
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 62.351132
$ws.Range("N8").Value = 187.053396
$ws.Range("O8").Value = 0.9620005273240912
$ws.Range("P8").Value = 0.9620005273240912
$ws.Range("Q8").Value = 12982.67636647072
$ws.Range("R8").Value = 116844.0872982365
$ws.Range("S8").Value = 0.4842230050783423
$ws.Range("T8").Value = 0.4842230050783423

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 208.2187756666666
$ws.Range("H9").Value = 624.6563269999999
$ws.Range("I9").Value = 0.5033500412159452
$ws.Range("J9").Value = 0.5033500412159452
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8562703333333334
$ws.Range("N9").Value = 2.568811
$ws.Range("O9").Value = 0.01321118776478095
$ws.Range("P9").Value = 0.01321118776478095
$ws.Range("Q9").Value = 178.2915604463552
$ws.Range("R9").Value = 1604.624044017197
$ws.Range("S9").Value = 0.006649851905914083
$ws.Range("T9").Value = 0.006649851905914083

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 208.2187756666666
$ws.Range("H10").Value = 624.6563269999999
$ws.Range("I10").Value = 0.5033500412159452
$ws.Range("J10").Value = 0.5033500412159452
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.606628666666667
$ws.Range("N10").Value = 4.819886
$ws.Range("O10").Value = 0.02478828491112776
$ws.Range("P10").Value = 0.02478828491112776
$ws.Range("Q10").Value = 334.5302539243024
$ws.Range("R10").Value = 3010.772285318722
$ws.Range("S10").Value = 0.01247718423168875
$ws.Range("T10").Value = 0.01247718423168875
